$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1338
$ws1.Range("F6").Value = 293
$ws1.Range("F7").Value = 1019
$ws1.Range("F8").Value = 10494
$ws1.Range("F9").Value = 12
$ws1.Range("F12").Value = 1028
$ws1.Range("F13").Value = 677
$ws1.Range("F14").Value = 11970
$ws1.Range("F15").Value = 12382
$ws1.Range("F17").Value = 113
$ws1.Range("F21").Value = 40

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 10

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 10
$ws4.Range("F6").Value = 1338
$ws4.Range("F7").Value = 293
$ws4.Range("F8").Value = 1019
$ws4.Range("F9").Value = 10494
$ws4.Range("F10").Value = 12
$ws4.Range("F13").Value = 1028
$ws4.Range("F14").Value = 677
$ws4.Range("F15").Value = 11970
$ws4.Range("F16").Value = 12382
$ws4.Range("F18").Value = 113
$ws4.Range("F22").Value = 40
